$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 414.14285
$ws.Range("J43").Value = 439.8
$ws.Range("L43").Value = 439.8
$ws.Range("N43").Value = -577.8
$ws.Range("H98").Value = 732.3871
$ws.Range("I98").Value = 790.5
$ws.Range("K98").Value = 790.5
$ws.Range("M98").Value = 707.5
$ws.Range("H99").Value = 1082
$ws.Range("I99").Value = 164
$ws.Range("K99").Value = 492
$ws.Range("M99").Value = 1006
$ws.Range("H118").Value = 739.25
$ws.Range("J118").Value = 777
$ws.Range("L118").Value = 2331
$ws.Range("N118").Value = -5645
$ws.Range("H122").Value = 732.3871
$ws.Range("I122").Value = 790.5
$ws.Range("K122").Value = 2371.5
$ws.Range("M122").Value = 78.5
$ws.Range("H123").Value = 39000
$ws.Range("J123").Value = 39000
$ws.Range("L123").Value = 39000
$ws.Range("N123").Value = -48800
$ws.Range("H127").Value = 1206.5385
$ws.Range("I127").Value = 881.3333
$ws.Range("J127").Value = 1485.2858
$ws.Range("K127").Value = 2643.9999
$ws.Range("L127").Value = 4455.857400000001
$ws.Range("M127").Value = 2316.0001
$ws.Range("N127").Value = -14375.8574
$ws.Range("H129").Value = 162746.6
$ws.Range("J129").Value = 165406.38
$ws.Range("L129").Value = 496219.14
$ws.Range("N129").Value = -506219.14
$ws.Range("H132").Value = 3588.5386
$ws.Range("I132").Value = 3926.348
$ws.Range("J132").Value = 998.6667
$ws.Range("K132").Value = 11779.044
$ws.Range("L132").Value = 2996.0001
$ws.Range("M132").Value = -9249.044
$ws.Range("N132").Value = -8056.0001
$ws.Range("H133").Value = 50704
$ws.Range("J133").Value = 50704
$ws.Range("L133").Value = 50704
$ws.Range("N133").Value = -60824
$ws.Range("H135").Value = 13518250
$ws.Range("I135").Value = 457.96774
$ws.Range("J135").Value = 83360180
$ws.Range("K135").Value = 4121.70966
$ws.Range("L135").Value = 750241620
$ws.Range("M135").Value = -1586.70966
$ws.Range("N135").Value = -750246690
$ws.Range("H137").Value = 35988.07
$ws.Range("I137").Value = 1479.95
$ws.Range("J137").Value = 112672.78
$ws.Range("K137").Value = 4439.85
$ws.Range("L137").Value = 338018.34
$ws.Range("M137").Value = -1889.85
$ws.Range("N137").Value = -343118.34
$ws.Range("H138").Value = 2004.2059
$ws.Range("I138").Value = 907.2308
$ws.Range("J138").Value = 2263.491
$ws.Range("K138").Value = 2721.6924
$ws.Range("L138").Value = 6790.473
$ws.Range("M138").Value = 2418.3076
$ws.Range("N138").Value = -17070.473

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H9").Value = 10009
$ws.Range("J9").Value = 10009
$ws.Range("L9").Value = 10009
$ws.Range("N9").Value = -10349
$ws.Range("H20").Value = 10009
$ws.Range("J20").Value = 10009
$ws.Range("L20").Value = 10009
$ws.Range("N20").Value = -10549
$ws.Range("H30").Value = 0
$ws.Range("J30").Value = 0
$ws.Range("L30").Value = 0
$ws.Range("N30").ClearContents()
$ws.Range("H63").Value = 3475377.5
$ws.Range("J63").Value = 5212166.5
$ws.Range("L63").Value = 5212166.5
$ws.Range("N63").Value = -5213538.5
$ws.Range("H66").Value = 3475377.5
$ws.Range("J66").Value = 5212166.5
$ws.Range("L66").Value = 26060832.5
$ws.Range("N66").Value = -26067696.5
$ws.Range("H74").Value = 41667824
$ws.Range("I74").Value = 55556150
$ws.Range("J74").Value = 2845.8333
$ws.Range("K74").Value = 55556150
$ws.Range("L74").Value = 2845.8333
$ws.Range("M74").Value = -55555276
$ws.Range("N74").Value = -4593.8333
$ws.Range("H77").Value = 41667824
$ws.Range("I77").Value = 55556150
$ws.Range("J77").Value = 2845.8333
$ws.Range("K77").Value = 277780750
$ws.Range("L77").Value = 14229.1665
$ws.Range("M77").Value = -277776382
$ws.Range("N77").Value = -22965.1665
$ws.Range("H122").Value = 3666.5
$ws.Range("I122").Value = 2400
$ws.Range("J122").Value = 9999
$ws.Range("K122").Value = 7200
$ws.Range("L122").Value = 29997
$ws.Range("M122").Value = -4750
$ws.Range("N122").Value = -34897
$ws.Range("H132").Value = 8367.958000000001
$ws.Range("I132").Value = 1226.7457
$ws.Range("K132").Value = 3680.2371
$ws.Range("M132").Value = -1150.2371

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 49420
$ws.Range("J138").Value = 49420
$ws.Range("L138").Value = 49420
$ws.Range("N138").Value = -59700

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10305.6
$ws.Range("I31").Value = 20210.295
$ws.Range("J31").Value = 2984.739
$ws.Range("K31").Value = 20210.295
$ws.Range("L31").Value = 2984.739
$ws.Range("M31").Value = -19915.295
$ws.Range("N31").Value = -3574.739
$ws.Range("H34").Value = 10305.6
$ws.Range("I34").Value = 20210.295
$ws.Range("J34").Value = 2984.739
$ws.Range("K34").Value = 20210.295
$ws.Range("L34").Value = 2984.739
$ws.Range("M34").Value = -20008.295
$ws.Range("N34").Value = -3388.739
$ws.Range("H52").Value = 39100
$ws.Range("J52").Value = 39100
$ws.Range("L52").Value = 39100
$ws.Range("N52").Value = -39688
$ws.Range("H62").Value = 111115224
$ws.Range("I62").Value = 200004000
$ws.Range("J62").Value = 4250.5
$ws.Range("K62").Value = 200004000
$ws.Range("L62").Value = 4250.5
$ws.Range("M62").Value = -200003376
$ws.Range("N62").Value = -5498.5
$ws.Range("H65").Value = 111115224
$ws.Range("I65").Value = 200004000
$ws.Range("J65").Value = 4250.5
$ws.Range("K65").Value = 1000020000
$ws.Range("L65").Value = 21252.5
$ws.Range("M65").Value = -1000016880
$ws.Range("N65").Value = -27492.5
$ws.Range("H132").Value = 16766.97
$ws.Range("I132").Value = 19589.678
$ws.Range("K132").Value = 58769.034
$ws.Range("M132").Value = -56239.034
$ws.Range("H134").Value = 847.5
$ws.Range("I134").Value = 767
$ws.Range("K134").Value = 2301
$ws.Range("M134").Value = 234

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H18").Value = 1058.3
$ws.Range("J18").Value = 4750
$ws.Range("L18").Value = 14250
$ws.Range("N18").Value = -14588
$ws.Range("H131").Value = 791.5859
$ws.Range("J131").Value = 816.70526
$ws.Range("L131").Value = 2450.11578
$ws.Range("N131").Value = -12530.11578

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 2412096.5
$ws.Range("I70").Value = 10206.8125
$ws.Range("K70").Value = 10206.8125
$ws.Range("M70").Value = -9936.8125
$ws.Range("H73").Value = 2412096.5
$ws.Range("I73").Value = 10206.8125
$ws.Range("K73").Value = 10206.8125
$ws.Range("M73").Value = -9270.8125
$ws.Range("H97").Value = 515.5
$ws.Range("J97").Value = 149.5
$ws.Range("L97").Value = 149.5
$ws.Range("N97").Value = -1141.5
$ws.Range("H132").Value = 66441.375
$ws.Range("I132").Value = 76513.21000000001
$ws.Range("J132").Value = 52340.8
$ws.Range("K132").Value = 229539.63
$ws.Range("L132").Value = 157022.4
$ws.Range("M132").Value = -227009.63
$ws.Range("N132").Value = -162082.4

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4852.52
$ws.Range("I61").Value = 3169.9333
$ws.Range("J61").Value = 7376.4
$ws.Range("K61").Value = 3169.9333
$ws.Range("L61").Value = 7376.4
$ws.Range("M61").Value = -2967.9333
$ws.Range("N61").Value = -7780.4
$ws.Range("H82").Value = 2771.5557
$ws.Range("I82").Value = 2420.5715
$ws.Range("K82").Value = 2420.5715
$ws.Range("M82").Value = -2059.5715
$ws.Range("H85").Value = 2771.5557
$ws.Range("I85").Value = 2420.5715
$ws.Range("K85").Value = 2420.5715
$ws.Range("M85").Value = -1172.5715
$ws.Range("H113").Value = 4852.52
$ws.Range("I113").Value = 3169.9333
$ws.Range("J113").Value = 7376.4
$ws.Range("K113").Value = 3169.9333
$ws.Range("L113").Value = 7376.4
$ws.Range("M113").Value = -999.9333000000001
$ws.Range("N113").Value = -11716.4
$ws.Range("H122").Value = 1636895.4
$ws.Range("I122").Value = 2803920.5
$ws.Range("J122").Value = 3060
$ws.Range("K122").Value = 8411761.5
$ws.Range("L122").Value = 9180
$ws.Range("M122").Value = -8409311.5
$ws.Range("N122").Value = -14080
$ws.Range("H132").Value = 2912.3635
$ws.Range("I132").Value = 2147.6
$ws.Range("J132").Value = 3549.6667
$ws.Range("K132").Value = 6442.799999999999
$ws.Range("L132").Value = 10649.0001
$ws.Range("M132").Value = -3912.799999999999
$ws.Range("N132").Value = -15709.0001

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()
$ws.Range("H27").Value = 35372.4
$ws.Range("J27").Value = 35372.4
$ws.Range("L27").Value = 35372.4
$ws.Range("N27").Value = -35510.4
$ws.Range("H132").Value = 2073.125
$ws.Range("I132").Value = 1161.4
$ws.Range("J132").Value = 3592.6667
$ws.Range("K132").Value = 3484.2
$ws.Range("L132").Value = 10778.0001
$ws.Range("M132").Value = -954.2000000000003
$ws.Range("N132").Value = -15838.0001
